$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-07-04 Thursday" "2024-07-05 Friday"

Replace-Text "985÷3=" "321÷6="
Replace-Text "632÷8=" "212÷6="
Replace-Text "874÷5=" "352÷2="
Replace-Text "570÷7=" "358÷6="
Replace-Text "945÷8=" "571÷5="
Replace-Text "732÷2=" "169÷6="
Replace-Text "324÷7=" "474÷2="
Replace-Text "494÷5=" "250÷7="
Replace-Text "873÷8=" "848÷2="
Replace-Text "832÷3=" "962÷6="
Replace-Text "632÷4=" "753÷4="
Replace-Text "311÷5=" "812÷4="
Replace-Text "810÷5=" "191÷5="
Replace-Text "565÷9=" "313÷7="
Replace-Text "786÷3=" "109÷3="
Replace-Text "705÷8=" "648÷3="
Replace-Text "687÷8=" "873÷7="
Replace-Text "875÷6=" "124÷3="
Replace-Text "477÷6=" "374÷8="
Replace-Text "833÷3=" "725÷5="
Replace-Text "196÷2=" "218÷6="
Replace-Text "541÷5=" "573÷7="
Replace-Text "803÷2=" "392÷5="
Replace-Text "207÷5=" "213÷3="
Replace-Text "370÷9=" "830÷7="
